# Orders for p7 and p8
# - Add a new "ITI" column (D) with values for trials 1-16.
# - Update ConditionType (column C) values for trials 1-16.
# - Remove trials 17-19 (rows 18-20), shrinking the table to 16 trials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "ITI"

# Updated ConditionType (C) and new ITI (D) values for trials 1-16 (rows 2-17)
$conditionType = @(2, 3, 2, 3, 3, 1, 2, 1, 4, 2, 4, 3, 4, 4, 1, 1)
$iti = @(6, 9, 8, 6, 6, 8, 8, 9, 6, 6, 7, 6, 8, 6, 7, 6)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $conditionType[$i]
    $ws.Cells.Item($row, 4).Value = $iti[$i]
}

# Remove the old trials 17-19 (rows 18-20), which no longer exist in the data
$ws.Range("A18:D20").Delete()

# Match the author's final selection state
$ws.Range("C18").Select() | Out-Null
